$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "69.376.42"
$ws.Range("E2").Value = "  +1.68%  "

# Row 3
$ws.Range("D3").Value = "3.381.64"
$ws.Range("E3").Value = "  +1.11%  "

# Row 4
$ws.Range("E4").Value = "  +0.13%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "580.96"
$ws.Range("E5").Value = "  -0.33%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "178.76"
$ws.Range("E6").Value = "  +1.02%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.00"
$ws.Range("E7").Value = "  +0.05%  "

# Row 8
$ws.Range("E8").Value = "  +0.60%  "

# Row 9
$ws.Range("E9").Value = "  +8.39%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.588"
$ws.Range("E10").Value = "  +1.00%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "48.49"
$ws.Range("E11").Value = "  +0.99%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0000284"
$ws.Range("E12").Value = "  +4.11%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "686.53"
$ws.Range("E13").Value = "  -0.74%  "

# Row 14
$ws.Range("E14").Value = "  +2.14%  "

# Row 15
$ws.Range("D15").Value = "3.925.76"
$ws.Range("E15").Value = "  +1.19%  "

# Row 16
$ws.Range("D16").Value = "69.498.97"
$ws.Range("E16").Value = "  +1.80%  "

# Row 17
$ws.Range("E17").Value = "  +0.81%  "

# Row 18
$ws.Range("D18").Value = "3.382.34"
$ws.Range("E18").Value = "  +1.69%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "17.85"
$ws.Range("E19").Value = "  +2.35%  "

# Row 20
$ws.Range("E20").Value = "  +0.91%  "

# Row 21
$ws.Range("E21").Value = "  +1.76%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "17.21"
$ws.Range("E22").Value = "  +1.46%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.35"
$ws.Range("E23").Value = "  -1.98%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "101.45"
$ws.Range("E24").Value = "  +1.55%  "

# Row 25
$ws.Range("E25").Value = "  -0.57%  "

# Row 26
$ws.Range("E26").Value = "  +0.12%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.74"
$ws.Range("E27").Value = "  +2.39%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "33.47"
$ws.Range("E28").Value = "  +1.49%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "8.73"
$ws.Range("E29").Value = "  +2.90%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "6.95"

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.80"
$ws.Range("E31").Value = "  +15.96%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "11.06"
$ws.Range("E32").Value = "  +0.19%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "551.13"
$ws.Range("E33").Value = "  -1.12%  "

# Row 34
$ws.Range("E34").Value = "  +0.22%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "57.79"
$ws.Range("E35").Value = "  +0.35%  "

# Row 36
$ws.Range("E36").Value = "  +0.03%  "

# Row 37
$ws.Range("D37").Value = "3.606.90"
$ws.Range("E37").Value = "  -2.42%  "

# Row 38
$ws.Range("E38").Value = "  +3.40%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "35.27"
$ws.Range("E39").Value = "  +1.54%  "

# Row 40
$ws.Range("D40").Value = "0.0₃0729"
$ws.Range("E40").Value = "  +8.60%  "

# Row 41
$ws.Range("E41").Value = "  +4.42%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.71"
$ws.Range("E42").Value = "  +3.90%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "3.39"
$ws.Range("E43").Value = "  +3.86%  "

# Row 44
$ws.Range("E44").Value = "  +3.40%  "

# Row 45
$ws.Range("E45").Value = "  +0.30%  "

# Row 46
$ws.Range("E46").Value = "  +0.41%  "

# Row 47
$ws.Range("E47").Value = "  +0.39%  "

# Row 48
$ws.Range("B48").Value = "Mantle"
$ws.Range("C48").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.39"
$ws.Range("E48").Value = "  +3.72%  "

# Row 49
$ws.Range("B49").Value = "FirstDigitalUSD"
$ws.Range("C49").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.00"
$ws.Range("E49").Value = "  -0.14%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "129.35"
$ws.Range("E50").Value = "  -1.12%  "

# Row 51
$ws.Range("E51").Value = "  +0.50%  "
